$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1961165048543689
$ws.Range("C2").Value = 0.5572815533980583
$ws.Range("J2").Value = 0.007766990291262136
$ws.Range("P2").Value = 0.1553398058252427
$ws.Range("S2").Value = 0.08349514563106795
$ws.Range("B3").Value = 0.006779661016949152
$ws.Range("C3").Value = 0.02033898305084746
$ws.Range("J3").Value = 0.03728813559322034
$ws.Range("P3").Value = 0.7491525423728813
$ws.Range("S3").Value = 0.1864406779661017
$ws.Range("J4").Value = 0.108433734939759
$ws.Range("P4").Value = 0.6987951807228916
$ws.Range("S4").Value = 0.1927710843373494
$ws.Range("J5").Value = 0.1666666666666667
$ws.Range("P5").Value = 0.6666666666666666
$ws.Range("S5").Value = 0.1666666666666667
$ws.Range("B6").Value = 0.05808080808080808
$ws.Range("D6").Value = 0.005050505050505051
$ws.Range("F6").Value = 0.07828282828282829
$ws.Range("J6").Value = 0.2297979797979798
$ws.Range("O6").Value = 0.02525252525252525
$ws.Range("Q6").Value = 0.1313131313131313
$ws.Range("R6").Value = 0.1035353535353535
$ws.Range("S6").Value = 0.3686868686868687
$ws.Range("B7").Value = 0.1428571428571428
$ws.Range("D7").Value = 0.008403361344537815
$ws.Range("E7").Value = 0.005602240896358543
$ws.Range("F7").Value = 0.06722689075630252
$ws.Range("J7").Value = 0.09523809523809523
$ws.Range("O7").Value = 0.02240896358543417
$ws.Range("Q7").Value = 0.1764705882352941
$ws.Range("R7").Value = 0.08123249299719888
$ws.Range("S7").Value = 0.4005602240896359
$ws.Range("B8").Value = 0.06406080347448426
$ws.Range("D8").Value = 0.01302931596091205
$ws.Range("E8").Value = 0.001085776330076004
$ws.Range("F8").Value = 0.0738327904451683
$ws.Range("J8").Value = 0.1205211726384365
$ws.Range("O8").Value = 0.01411509229098806
$ws.Range("Q8").Value = 0.1672095548317047
$ws.Range("R8").Value = 0.1107491856677524
$ws.Range("S8").Value = 0.4353963083604778
$ws.Range("B9").Value = 0.05955334987593052
$ws.Range("D9").Value = 0.02233250620347394
$ws.Range("F9").Value = 0.04962779156327544
$ws.Range("J9").Value = 0.1240694789081886
$ws.Range("O9").Value = 0.01985111662531017
$ws.Range("Q9").Value = 0.1811414392059553
$ws.Range("R9").Value = 0.1389578163771712
$ws.Range("S9").Value = 0.4044665012406948
$ws.Range("B10").Value = 0.105708245243129
$ws.Range("D10").Value = 0.0241014799154334
$ws.Range("E10").Value = 0.002114164904862579
$ws.Range("F10").Value = 0.06004228329809725
$ws.Range("J10").Value = 0.1154334038054968
$ws.Range("O10").Value = 0.01649048625792812
$ws.Range("Q10").Value = 0.2029598308668076
$ws.Range("R10").Value = 0.09852008456659619
$ws.Range("S10").Value = 0.3746300211416491
$ws.Range("G11").Value = 0.131578947368421
$ws.Range("J11").Value = 0.07330827067669173
$ws.Range("K11").Value = 0.1766917293233083
$ws.Range("L11").Value = 0.6071428571428571
$ws.Range("S11").Value = 0.0112781954887218
$ws.Range("G12").Value = 0.7099697885196374
$ws.Range("J12").Value = 0.2386706948640483
$ws.Range("K12").Value = 0.00906344410876133
$ws.Range("L12").Value = 0.01812688821752266
$ws.Range("S12").Value = 0.02416918429003021
$ws.Range("G13").Value = 0.7065217391304348
$ws.Range("J13").Value = 0.2608695652173913
$ws.Range("S13").Value = 0.03260869565217391
$ws.Range("F15").Value = 0.01121076233183856
$ws.Range("H15").Value = 0.1569506726457399
$ws.Range("I15").Value = 0.05605381165919283
$ws.Range("J15").Value = 0.3542600896860987
$ws.Range("K15").Value = 0.08520179372197309
$ws.Range("M15").Value = 0.02017937219730942
$ws.Range("O15").Value = 0.06502242152466367
$ws.Range("S15").Value = 0.2511210762331839
$ws.Range("F16").Value = 0.01704545454545454
$ws.Range("H16").Value = 0.1590909090909091
$ws.Range("I16").Value = 0.08522727272727272
$ws.Range("J16").Value = 0.4176136363636364
$ws.Range("K16").Value = 0.1107954545454545
$ws.Range("N16").Value = 0.005681818181818182
$ws.Range("O16").Value = 0.05397727272727273
$ws.Range("S16").Value = 0.1278409090909091
$ws.Range("F17").Value = 0.009864364981504316
$ws.Range("H17").Value = 0.1911220715166461
$ws.Range("I17").Value = 0.09124537607891492
$ws.Range("J17").Value = 0.4204685573366215
$ws.Range("K17").Value = 0.09247842170160296
$ws.Range("M17").Value = 0.01479654747225647
$ws.Range("N17").Value = 0.001233045622688039
$ws.Range("O17").Value = 0.07891491985203453
$ws.Range("S17").Value = 0.0998766954377312
$ws.Range("F18").Value = 0.01538461538461539
$ws.Range("H18").Value = 0.2263736263736264
$ws.Range("I18").Value = 0.08571428571428572
$ws.Range("J18").Value = 0.410989010989011
$ws.Range("K18").Value = 0.08571428571428572
$ws.Range("M18").Value = 0.01758241758241758
$ws.Range("N18").Value = 0.002197802197802198
$ws.Range("O18").Value = 0.05494505494505494
$ws.Range("S18").Value = 0.1010989010989011
$ws.Range("F19").Value = 0.01924651924651925
$ws.Range("H19").Value = 0.2178542178542179
$ws.Range("I19").Value = 0.09623259623259624
$ws.Range("J19").Value = 0.3484848484848485
$ws.Range("K19").Value = 0.0990990990990991
$ws.Range("M19").Value = 0.02375102375102375
$ws.Range("O19").Value = 0.06920556920556921
$ws.Range("S19").Value = 0.1261261261261261
